# Applies "done" markers to the UC sheet, in a new column G, for rows
# corresponding to use cases that have been completed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("UC")

# Rows whose use case is now marked as "done"
$doneRows = @(2, 3, 10, 12, 13, 14, 23, 29)

foreach ($r in $doneRows) {
    $ws.Cells.Item($r, 7).Value = "done"
}

# Activate the sheet and update the view/selection to match the edit
$ws.Activate()
$ws.Application.ActiveWindow.ScrollRow = 22
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("G2").Select()
